$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the data block (row 179),
# pushing all subsequent rows (old 179-205) down by one.
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row 179 with the new week's values.
$ws.Range("A179").Value = 3
$ws.Range("B179").Value = "Femacal de La Calera"
$ws.Range("C179").Value = "Coquimbo"
$ws.Range("D179").Value = 44474
$ws.Range("E179").Value = 5
$ws.Range("F179").Value = 100112043
$ws.Range("G179").Value = "Pepino ensalada"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 80
$ws.Range("K179").Value = 15500
$ws.Range("L179").Value = 16000
$ws.Range("M179").Value = 15750
$ws.Range("N179").Value = "$/caja 70 unidades"
$ws.Range("O179").Value = "Región de Arica y Parinacota"
$ws.Range("P179").Value = 225
$ws.Range("Q179").Value = 70
$ws.Range("R179").Value = "Hortaliza"
